$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Images for Upload")

# Fill in missing required values before removing the image-category column.
# Row 5 ("Head, Tail, Fore Edge") is Binding evidence, like rows 2-4.
$ws.Range("H5").Value = "Binding"
# Row 9 (Title Page) is a non-evidence image, flagged explicitly.
$ws.Range("H9").Value = "Title Page (non-evidence)"

# Remove the "evidence: image category" column entirely (image_type removal).
$ws.Columns.Item(7).Select()
$ws.Columns.Item(7).Delete()

$wb.Save()
